# Update cryptocurrency price (column D) and 1h volume change (column E)
# for each coin row to reflect the latest refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPrices = @{
    2 = "28.245.53"
    3 = "1.814.93"
    5 = "325.89"
    7 = "0.4361"
    8 = "0.3668"
    9 = "44.89"
    10 = "0.07684"
    11 = "1.141"
    12 = "1.001"
    13 = "22.01"
    14 = "6.316"
    15 = "7.497"
    16 = "1.819.48"
    17 = "95.07"
    18 = "0.00001082"
    19 = "0.06483"
    20 = "0.9997"
    22 = "6.238"
    23 = "28.274.24"
    24 = "11.59"
    25 = "2.131"
    26 = "161.26"
    27 = "20.72"
    28 = "2.025.55"
    29 = "2.283"
    30 = "130.10"
    31 = "1.214"
    32 = "6.028"
    33 = "0.09144"
    35 = "13.12"
    36 = "0.02382"
    37 = "5.251"
    39 = "0.6604"
    42 = "8.070"
    43 = "1.429"
    44 = "0.9991"
    45 = "13.85"
    46 = "0.6115"
    47 = "3.741"
    48 = "2.025"
    49 = "125.83"
    50 = "1.164"
    51 = "0.06996"
}

$newVolumes = @{
    2 = "  +2.09%  "
    3 = "  +3.59%  "
    4 = "  -0.27%  "
    5 = "  +0.50%  "
    6 = "  +0.01%  "
    7 = "  +0.91%  "
    8 = "  +0.49%  "
    9 = "  -1.22%  "
    10 = "  +3.06%  "
    11 = "  +2.04%  "
    12 = "  -0.01%  "
    13 = "  +1.97%  "
    14 = "  +2.60%  "
    15 = "  +3.53%  "
    16 = "  +3.91%  "
    17 = "  +7.73%  "
    18 = "  +1.48%  "
    19 = "  +4.52%  "
    20 = "  -0.01%  "
    21 = "  +1.82%  "
    22 = "  +1.53%  "
    23 = "  +2.11%  "
    24 = "  -0.39%  "
    25 = "  -8.43%  "
    26 = "  +4.96%  "
    27 = "  +0.73%  "
    28 = "  +3.90%  "
    29 = "  -3.33%  "
    30 = "  +1.95%  "
    31 = "  -0.31%  "
    32 = "  +5.43%  "
    33 = "  -0.13%  "
    34 = "  -1.98%  "
    35 = "  +4.02%  "
    36 = "  +3.26%  "
    37 = "  +3.02%  "
    38 = "  +0.85%  "
    39 = "  +2.18%  "
    40 = "  +1.95%  "
    41 = "  +0.29%  "
    42 = "  +1.67%  "
    43 = "  +0.15%  "
    44 = "  +0.01%  "
    45 = "  +0.40%  "
    46 = "  +3.14%  "
    47 = "  -0.15%  "
    48 = "  +2.70%  "
    49 = "  +0.06%  "
    50 = "  +3.90%  "
    51 = "  +1.40%  "
}

foreach ($row in $newPrices.Keys) {
    $priceCell = $ws.Range("D$row")
    $newPrice = $newPrices[$row]
    # Force text storage so values like "44.89" or "8.070" keep their exact
    # printed form instead of being reinterpreted as numbers.
    if ($newPrice -match "^-?\d+(\.\d+)?$") {
        $priceCell.NumberFormat = "@"
    }
    $priceCell.Value = $newPrice
}

foreach ($row in $newVolumes.Keys) {
    $ws.Range("E$row").Value = $newVolumes[$row]
}
